# Refresh cryptos list with the latest scraped Price / Volume(1h) figures,
# and pick up the current CoinRanking rank reshuffle for a few coins near the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coins that changed rank position this run: full row (Coin/Link/Price/Volume) refresh.
$reorderedCoins = @(
    @{ Row = 44; Coin = 'Bittensor'; Link = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; Price = '284.11'; Volume = '  -6.89%  ' },
    @{ Row = 45; Coin = 'InjectiveProtocol'; Link = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; Price = '20.88'; Volume = '  -6.33%  ' },
    @{ Row = 46; Coin = 'FirstDigitalUSD'; Link = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; Price = '0.999'; Volume = '  -0.06%  ' },
    @{ Row = 48; Coin = 'dogwifhat'; Link = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; Price = '1.92'; Volume = '  -11.66%  ' },
    @{ Row = 49; Coin = 'WhiteBITCoin'; Link = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'; Price = '10.44'; Volume = '  +0.61%  ' }
)

foreach ($u in $reorderedCoins) {
    $r = $u.Row
    $ws.Range("B$r").Value = $u.Coin
    $ws.Range("C$r").Value = $u.Link
    # Force Price to stay text (avoid Excel auto-converting "284.11" style strings to numbers).
    $priceCell = $ws.Range("D$r")
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $u.Price
    $priceCell.Style = "Normal"
    $ws.Range("E$r").Value = $u.Volume
}

# Remaining coins keep their position; only Price and/or Volume(1h) moved.
$priceVolumeUpdates = @(
    @{ Row = 2; Price = '64.125.78'; Volume = '  -0.81%  ' },
    @{ Row = 3; Price = '3.145.56' },
    @{ Row = 4; Price = '1.00'; Volume = '  +0.01%  ' },
    @{ Row = 5; Price = '567.08'; Volume = '  -1.86%  ' },
    @{ Row = 6; Price = '163.59'; Volume = '  -6.65%  ' },
    @{ Row = 7; Volume = '  -0.01%  ' },
    @{ Row = 8; Price = '0.585'; Volume = '  -7.06%  ' },
    @{ Row = 9; Price = '6.69'; Volume = '  -0.99%  ' },
    @{ Row = 10; Price = '0.116'; Volume = '  -5.39%  ' },
    @{ Row = 11; Price = '0.380'; Volume = '  -2.80%  ' },
    @{ Row = 12; Price = '3.691.23'; Volume = '  -2.86%  ' },
    @{ Row = 13; Volume = '  -0.83%  ' },
    @{ Row = 14; Price = '63.885.34'; Volume = '  -1.43%  ' },
    @{ Row = 15; Price = '24.98'; Volume = '  -2.93%  ' },
    @{ Row = 16; Price = '3.132.02'; Volume = '  -3.52%  ' },
    @{ Row = 17; Price = '0.0000154'; Volume = '  -3.11%  ' },
    @{ Row = 18; Price = '409.84'; Volume = '  -1.26%  ' },
    @{ Row = 19; Price = '5.23'; Volume = '  -2.58%  ' },
    @{ Row = 20; Price = '12.48'; Volume = '  -2.88%  ' },
    @{ Row = 21; Price = '7.08'; Volume = '  -1.71%  ' },
    @{ Row = 22; Volume = '  +0.48%  ' },
    @{ Row = 23; Price = '67.56'; Volume = '  -3.75%  ' },
    @{ Row = 24; Price = '0.199'; Volume = '  -1.78%  ' },
    @{ Row = 25; Price = '0.481'; Volume = '  -3.16%  ' },
    @{ Row = 26; Volume = '  -7.02%  ' },
    @{ Row = 27; Price = '8.83'; Volume = '  -3.77%  ' },
    @{ Row = 28; Price = '0.997'; Volume = '  -0.36%  ' },
    @{ Row = 29; Volume = '  -3.98%  ' },
    @{ Row = 30; Price = '21.07'; Volume = '  -3.48%  ' },
    @{ Row = 31; Price = '6.25'; Volume = '  -3.05%  ' },
    @{ Row = 32; Price = '4.76'; Volume = '  -5.23%  ' },
    @{ Row = 33; Price = '1.12'; Volume = '  -3.22%  ' },
    @{ Row = 34; Price = '152.94'; Volume = '  -2.58%  ' },
    @{ Row = 35; Volume = '  -5.10%  ' },
    @{ Row = 36; Price = '2.747.94'; Volume = '  -2.90%  ' },
    @{ Row = 37; Price = '1.66'; Volume = '  -4.57%  ' },
    @{ Row = 38; Volume = '  -8.28%  ' },
    @{ Row = 39; Price = '4.06'; Volume = '  -3.57%  ' },
    @{ Row = 40; Volume = '  -4.29%  ' },
    @{ Row = 41; Price = '0.0625'; Volume = '  -0.24%  ' },
    @{ Row = 42; Price = '5.43'; Volume = '  -5.64%  ' },
    @{ Row = 43; Price = '0.0258'; Volume = '  -2.04%  ' },
    @{ Row = 47; Price = '0.0975'; Volume = '  -3.37%  ' },
    @{ Row = 50; Price = '5.69'; Volume = '  -2.25%  ' },
    @{ Row = 51; Price = '0.887'; Volume = '  -5.34%  ' }
)

foreach ($u in $priceVolumeUpdates) {
    $r = $u.Row
    if ($u.ContainsKey("Price")) {
        $priceCell = $ws.Range("D$r")
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.Price
        $priceCell.Style = "Normal"
    }
    if ($u.ContainsKey("Volume")) {
        $ws.Range("E$r").Value = $u.Volume
    }
}
